$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefixMap = @{
    2  = "B-"
    3  = "B-"
    4  = "B-"
    5  = "V-"
    6  = "B-"
    7  = "B-"
    8  = "B-"
    9  = "B-"
    10 = "B-"
    11 = "B-"
    12 = "B-"
    13 = "B-"
    14 = "B-"
    15 = "V-"
    16 = "V-"
    17 = "B-"
    18 = "B-"
}

foreach ($row in $prefixMap.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $current = $cell.Value2
    $cell.Value = $prefixMap[$row] + $current
}
